$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Paulo Rodrigues"
$ws.Range("B3").Value = "Academia2024foco"

$ws.Range("B8").Select()
